# Refresh the cryptos price/volume snapshot (GitHub Actions cron update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") ---
# Values are stored as plain text (e.g. "62.779.22", "1.00", "0.0000222"),
# several of which look numeric to Excels auto-detection and would silently
# be coerced into a Double (losing the exact formatting / becoming sci-notation)
# if assigned directly. Forcing the Text number format before the write keeps
# the literal string, and ClearFormats() afterwards drops that temporary format
# again so the cell ends up with the same (default/no) style as the source file.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.779.22"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.074.68"
$ws.Range("D3").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.40"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.50"
$ws.Range("D6").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.071.73"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.487"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.153"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.12"
$ws.Range("D11").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000222"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.36"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.539.59"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.848.94"
$ws.Range("D16").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.076.21"
$ws.Range("D18").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.84"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.26"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.697"
$ws.Range("D22").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.52"
$ws.Range("D24").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.21"
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.74"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.86"
$ws.Range("D31").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "57.90"
$ws.Range("D33").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.20"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "470.92"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.106.85"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0389"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0787"
$ws.Range("D40").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.57"
$ws.Range("D43").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.26"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.19"
$ws.Range("D48").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0509"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.99"
$ws.Range("D51").ClearFormats()

# --- Column E ("Volume(1h)") ---
# These are always padded percentage strings (e.g. "  -4.11%  "); the "%" plus
# surrounding spaces keep Excel from treating them as numbers, so a plain Value
# assignment is sufficient to preserve them as text verbatim.
$ws.Range("E2").Value = "  -4.11%  "
$ws.Range("E3").Value = "  -3.61%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -5.72%  "
$ws.Range("E6").Value = "  -9.81%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -3.42%  "
$ws.Range("E9").Value = "  -3.54%  "
$ws.Range("E10").Value = "  -4.08%  "
$ws.Range("E11").Value = "  -11.36%  "
$ws.Range("E12").Value = "  -4.74%  "
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("E14").Value = "  -9.73%  "
$ws.Range("E15").Value = "  -4.48%  "
$ws.Range("E16").Value = "  -4.16%  "
$ws.Range("E17").Value = "  -2.96%  "
$ws.Range("E18").Value = "  -3.72%  "
$ws.Range("E19").Value = "  -6.48%  "
$ws.Range("E20").Value = "  -9.76%  "
$ws.Range("E21").Value = "  -7.77%  "
$ws.Range("E22").Value = "  -4.77%  "
$ws.Range("E23").Value = "  -6.65%  "
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("E25").Value = "  -9.23%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -7.20%  "
$ws.Range("E28").Value = "  -11.15%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  -5.04%  "
$ws.Range("E31").Value = "  -16.39%  "
$ws.Range("E32").Value = "  -5.51%  "
$ws.Range("E33").Value = "  +6.15%  "
$ws.Range("E34").Value = "  -10.91%  "
$ws.Range("E35").Value = "  -5.23%  "
$ws.Range("E36").Value = "  -6.20%  "
$ws.Range("E37").Value = "  -15.15%  "
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("E39").Value = "  -12.26%  "
$ws.Range("E40").Value = "  -6.69%  "
$ws.Range("E41").Value = "  -5.56%  "
$ws.Range("E42").Value = "  -11.29%  "
$ws.Range("E43").Value = "  -10.12%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  -9.78%  "
$ws.Range("E46").Value = "  -11.97%  "
$ws.Range("E47").Value = "  -7.74%  "
$ws.Range("E48").Value = "  -4.37%  "
$ws.Range("E49").Value = "  -3.59%  "
$ws.Range("E50").Value = "  -5.83%  "
$ws.Range("E51").Value = "  -8.45%  "
